$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 228, shifting rows 228:314 down to 229:315.
$ws.Rows.Item(228).Insert()

# Populate the newly-inserted row 228 with its data.
$ws.Range("A228").Value = 9
$ws.Range("B228").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C228").Value = "Metropolitana"
$ws.Range("D228").Value = 44988
$ws.Range("E228").Value = 13
$ws.Range("F228").Value = "Fruta"
$ws.Range("G228").Value = 100101
$ws.Range("H228").Value = "Berries"
$ws.Range("I228").Value = 100101001
$ws.Range("J228").Value = "Arándano (blue)"
$ws.Range("K228").Value = "Sin especificar"
$ws.Range("L228").Value = "Primera"
$ws.Range("M228").Value = 300
$ws.Range("N228").Value = 3000
$ws.Range("O228").Value = 3000
$ws.Range("P228").Value = 3000
$ws.Range("Q228").Value = "$/bandeja 2 kilos"
$ws.Range("R228").Value = "Provincia de Curicó"
$ws.Range("S228").Value = 1500
$ws.Range("T228").Value = 2
